$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# URL: matchsource -> matchsync
$ws1.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-abo-codes"

# Experimental: blank -> true (must be stored as literal text "true", not a
# Boolean, so build it as a text formula result and paste-as-values to strip
# the formula while keeping the Text type / existing cell style)
$ws1.Range("B7").Formula = "=""true"""
$ws1.Range("B7").Copy()
$ws1.Range("B7").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Date: updated timestamp
$ws1.Range("B8").Value = "2024-02-19T18:37:26-06:00"
